$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''65.221.90'
$ws.Range('E2').Value = '  -4.63%  '
$ws.Range('D3').Value = '''3.379.38'
$ws.Range('E3').Value = '  -6.57%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = '''187.62'
$ws.Range('E5').Value = '  -7.63%  '
$ws.Range('D6').Value = '''529.03'
$ws.Range('E6').Value = '  -6.91%  '
$ws.Range('D7').Value = '''0.603'
$ws.Range('E7').Value = '  -3.28%  '
$ws.Range('D8').Value = '''3.372.63'
$ws.Range('E8').Value = '  -6.66%  '
$ws.Range('E9').Value = '  +0.00%  '
$ws.Range('D10').Value = '''0.627'
$ws.Range('E10').Value = '  -7.41%  '
$ws.Range('D11').Value = '''58.50'
$ws.Range('E11').Value = '  -4.76%  '
$ws.Range('E12').Value = '  -12.08%  '
$ws.Range('D13').Value = '''0.0000254'
$ws.Range('E13').Value = '  -11.65%  '
$ws.Range('D14').Value = '''9.25'
$ws.Range('E14').Value = '  -8.02%  '
$ws.Range('D15').Value = '''3.911.35'
$ws.Range('E15').Value = '  -6.82%  '
$ws.Range('E16').Value = '  -3.05%  '
$ws.Range('D17').Value = '''3.377.92'
$ws.Range('E17').Value = '  -6.61%  '
$ws.Range('D18').Value = '''64.957.78'
$ws.Range('E18').Value = '  -4.78%  '
$ws.Range('D19').Value = '''17.46'
$ws.Range('E19').Value = '  -8.42%  '
$ws.Range('D20').Value = '''11.13'
$ws.Range('E20').Value = '  -10.27%  '
$ws.Range('D21').Value = '''0.973'
$ws.Range('E21').Value = '  -9.83%  '
$ws.Range('D22').Value = '''371.72'
$ws.Range('E22').Value = '  -8.05%  '
$ws.Range('D23').Value = '''81.48'
$ws.Range('E23').Value = '  -4.81%  '
$ws.Range('D24').Value = '''3.72'
$ws.Range('E24').Value = '  -11.12%  '
$ws.Range('D25').Value = '''10.78'
$ws.Range('E25').Value = '  -19.48%  '
$ws.Range('D26').Value = '''3.76'
$ws.Range('E26').Value = '  -4.48%  '
$ws.Range('B27').Value = 'LEO'
$ws.Range('C27').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D27').Value = '''5.81'
$ws.Range('E27').Value = '  -5.27%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''11.51'
$ws.Range('E28').Value = '  -8.79%  '
$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D29').Value = '''2.65'
$ws.Range('E29').Value = '  -10.39%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = '''8.53'
$ws.Range('E30').Value = '  -8.96%  '
$ws.Range('B31').Value = 'EthereumClassic'
$ws.Range('C31').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D31').Value = '''29.62'
$ws.Range('E31').Value = '  -6.44%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').Value = '''669.87'
$ws.Range('E32').Value = '  -0.14%  '
$ws.Range('B33').Value = 'NEARProtocol'
$ws.Range('C33').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D33').Value = '''6.72'
$ws.Range('E33').Value = '  -16.71%  '
$ws.Range('B34').Value = 'Cosmos'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D34').Value = '''11.18'
$ws.Range('E34').Value = '  -9.38%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''61.17'
$ws.Range('E35').Value = '  -4.11%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.106'
$ws.Range('E36').Value = '  -7.90%  '
$ws.Range('B37').Value = 'Dai'
$ws.Range('C37').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D37').Value = '''1.00'
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('B38').Value = 'InjectiveProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D38').Value = '''36.54'
$ws.Range('E38').Value = '  -13.65%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '''0.379'
$ws.Range('E39').Value = '  -10.31%  '
$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').Value = '''0.998'
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.127'
$ws.Range('E41').Value = '  -6.73%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '''2.854.57'
$ws.Range('E42').Value = '  -11.25%  '
$ws.Range('B43').Value = 'ThetaToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D43').Value = '''2.75'
$ws.Range('E43').Value = '  -14.33%  '
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''2.66'
$ws.Range('E44').Value = '  -7.82%  '
$ws.Range('B45').Value = 'VeChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D45').Value = '''0.0394'
$ws.Range('E45').Value = '  -6.10%  '
$ws.Range('B46').Value = 'PEPE'
$ws.Range('C46').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D46').Value = '''0.0₃0621'
$ws.Range('E46').Value = '  -20.61%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '''2.35'
$ws.Range('E47').Value = '  -15.55%  '
$ws.Range('B48').Value = 'Monero'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D48').Value = '''137.52'
$ws.Range('E48').Value = '  -1.27%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '''0.125'
$ws.Range('E49').Value = '  -5.46%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').Value = '''2.81'
$ws.Range('E50').Value = '  -8.69%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''2.58'
$ws.Range('E51').Value = '  -6.56%  '
